# SMILocalization.xlsx update:
#  - Add "Session Expired" / "Pilot session no longer valid" strings for the
#    new PILOT_MODE_EXPIRED_TITLE / PILOT_MODE_EXPIRED_DESCRIPTION keys.
#  - Rename the pilot title copy from "SMI PILOT" to "SMI Pilot".
#  - Update the pilot description copy from "...press Play to start" to
#    "...press Play to continue".
#  - Widen column A now that "PILOT_MODE_EXPIRED_DESCRIPTION" is the longest
#    key, and drop the stale best-fit flag.
#  - Leave the selection on B29, matching the author's last edit location.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows right above the existing pilot-mode rows (28 & 29)
# for the new "session expired" title/description pair.
$ws.Rows.Item(28).Resize(2).Insert()

# Give the two new "C" (nl) cells the same red "Bad" (awaiting translation)
# style already used by the rest of the pilot-mode block, by copying the
# format from the row directly below rather than re-applying the named
# style (which would otherwise mint a duplicate style entry).
$ws.Range("C30").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A28").Value = "PILOT_MODE_EXPIRED_TITLE"
$ws.Range("B28").Value = "Session Expired"
$ws.Range("C28").Value = "Session Expired"
$ws.Rows.Item(28).RowHeight = 14.5

$ws.Range("A29").Value = "PILOT_MODE_EXPIRED_DESCRIPTION"
$ws.Range("B29").Value = "Pilot session no longer valid"
$ws.Range("C29").Value = "Pilot session no longer valid"
$ws.Rows.Item(29).RowHeight = 14.5

# Former row 28 (PILOT_MODE_TITLE) is now row 30 - fix the casing of the copy.
$ws.Range("B30").Value = "SMI Pilot"
$ws.Range("C30").Value = "SMI Pilot"

# Former row 29 (PILOT_MODE_DESCRIPTION) is now row 31 - "start" -> "continue".
$ws.Range("B31").Value = "Welcome to the pilot, press Play to continue"
$ws.Range("C31").Value = "Welcome to the pilot, press Play to continue"

# Column A needs to be a bit wider to fit the new, longer keys; this also
# clears the old best-fit flag since the width is now explicit.
$ws.Columns.Item(1).ColumnWidth = 39

# Restore the cursor to where the author left it.
$ws.Range("B29").Select()
